$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "51.593.76"
Set-TextValue $ws.Range("E2") "  -1.34%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.775.76"
Set-TextValue $ws.Range("E3") "  -2.10%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "358.12"
Set-TextValue $ws.Range("E5") "  -0.72%  "

# Row 6
Set-TextValue $ws.Range("D6") "109.17"
Set-TextValue $ws.Range("E6") "  -3.09%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -3.50%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.08%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.586"
Set-TextValue $ws.Range("E9") "  -2.87%  "

# Row 10
Set-TextValue $ws.Range("D10") "39.70"
Set-TextValue $ws.Range("E10") "  -3.47%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.139"
Set-TextValue $ws.Range("E11") "  +4.63%  "

# Row 12
Set-TextValue $ws.Range("E12") "  -3.68%  "

# Row 13
Set-TextValue $ws.Range("D13") "19.67"
Set-TextValue $ws.Range("E13") "  -1.63%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.58"
Set-TextValue $ws.Range("E14") "  -2.95%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.214.16"
Set-TextValue $ws.Range("E15") "  -1.94%  "

# Row 16
Set-TextValue $ws.Range("D16") "2.797.07"
Set-TextValue $ws.Range("E16") "  -0.76%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.927"
Set-TextValue $ws.Range("E17") "  -0.45%  "

# Row 18
Set-TextValue $ws.Range("D18") "51.555.01"
Set-TextValue $ws.Range("E18") "  -1.21%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.58"
Set-TextValue $ws.Range("E19") "  +0.42%  "

# Row 20
Set-TextValue $ws.Range("D20") "3.09"
Set-TextValue $ws.Range("E20") "  -2.25%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.20"
Set-TextValue $ws.Range("E21") "  -2.41%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -3.32%  "

# Row 23
Set-TextValue $ws.Range("D23") "70.12"
Set-TextValue $ws.Range("E23") "  -0.83%  "

# Row 24
Set-TextValue $ws.Range("D24") "267.91"
Set-TextValue $ws.Range("E24") "  -1.97%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -3.24%  "

# Row 26
Set-TextValue $ws.Range("D26") "26.30"
Set-TextValue $ws.Range("E26") "  -2.67%  "

# Row 27
Set-TextValue $ws.Range("B27") "Kaspa"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.168"
Set-TextValue $ws.Range("E27") "  +17.11%  "

# Row 28
Set-TextValue $ws.Range("B28") "Dai"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D28") "1.00"
Set-TextValue $ws.Range("E28") "  -0.09%  "

# Row 29
Set-TextValue $ws.Range("D29") "10.17"
Set-TextValue $ws.Range("E29") "  -1.86%  "

# Row 30
Set-TextValue $ws.Range("D30") "2.19"
Set-TextValue $ws.Range("E30") "  -3.07%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.21"
Set-TextValue $ws.Range("E31") "  +4.55%  "

# Row 32
Set-TextValue $ws.Range("D32") "35.30"
Set-TextValue $ws.Range("E32") "  -0.11%  "

# Row 33
Set-TextValue $ws.Range("D33") "52.03"
Set-TextValue $ws.Range("E33") "  -0.22%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0443"
Set-TextValue $ws.Range("E34") "  -8.42%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -1.82%  "

# Row 36
Set-TextValue $ws.Range("D36") "5.23"
Set-TextValue $ws.Range("E36") "  -6.83%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.09%  "

# Row 38
Set-TextValue $ws.Range("D38") "18.69"
Set-TextValue $ws.Range("E38") "  +1.22%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.13"
Set-TextValue $ws.Range("E39") "  -4.99%  "

# Row 40
Set-TextValue $ws.Range("E40") "  -5.07%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -3.13%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.51"
Set-TextValue $ws.Range("E42") "  -1.28%  "

# Row 43
Set-TextValue $ws.Range("B43") "WEMIXToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "2.20"
Set-TextValue $ws.Range("E43") "  -4.26%  "

# Row 44
Set-TextValue $ws.Range("B44") "Monero"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D44") "119.52"
Set-TextValue $ws.Range("E44") "  -6.13%  "

# Row 45
Set-TextValue $ws.Range("D45") "21.76"
Set-TextValue $ws.Range("E45") "  -5.86%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.081.70"
Set-TextValue $ws.Range("E46") "  -0.79%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.25"
Set-TextValue $ws.Range("E47") "  -3.87%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.29"
Set-TextValue $ws.Range("E48") "  -0.51%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.940"
Set-TextValue $ws.Range("E49") "  -3.11%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -6.38%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.29"
Set-TextValue $ws.Range("E51") "  +5.06%  "
